$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($Range, $Text) {
    # Force the cell to remain plain text (column is entirely text-typed in
    # the source data, e.g. "72.420.13" / "0.170" / "1.00") so that Excel's
    # automatic number detection does not strip meaningful trailing zeros
    # or re-interpret multi-dot strings.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-CellText $ws.Range("D2") "72.389.79"
Set-CellText $ws.Range("E2") "  +1.62%  "
Set-CellText $ws.Range("D3") "2.684.37"
Set-CellText $ws.Range("E3") "  +1.86%  "
Set-CellText $ws.Range("E4") "  +0.04%  "
Set-CellText $ws.Range("D5") "600.72"
Set-CellText $ws.Range("E5") "  -1.09%  "
Set-CellText $ws.Range("D6") "177.42"
Set-CellText $ws.Range("E6") "  -2.13%  "
Set-CellText $ws.Range("D8") "0.524"
Set-CellText $ws.Range("E8") "  -0.04%  "
Set-CellText $ws.Range("D9") "2.682.94"
Set-CellText $ws.Range("E9") "  +1.89%  "
Set-CellText $ws.Range("D10") "0.170"
Set-CellText $ws.Range("E10") "  +2.64%  "
Set-CellText $ws.Range("E11") "  +2.24%  "
Set-CellText $ws.Range("E12") "  +2.28%  "
Set-CellText $ws.Range("D13") "5.03"
Set-CellText $ws.Range("E13") "  +0.20%  "
Set-CellText $ws.Range("D14") "3.176.62"
Set-CellText $ws.Range("E14") "  +4.10%  "
Set-CellText $ws.Range("E15") "  -1.31%  "
Set-CellText $ws.Range("D16") "72.302.25"
Set-CellText $ws.Range("E16") "  +1.64%  "
Set-CellText $ws.Range("D17") "26.35"
Set-CellText $ws.Range("E17") "  -1.31%  "
Set-CellText $ws.Range("D18") "2.686.41"
Set-CellText $ws.Range("E18") "  +2.12%  "
Set-CellText $ws.Range("D19") "12.05"
Set-CellText $ws.Range("E19") "  +4.37%  "
Set-CellText $ws.Range("E20") "  +0.14%  "
Set-CellText $ws.Range("D21") "371.39"
Set-CellText $ws.Range("E21") "  -2.97%  "
Set-CellText $ws.Range("E22") "  +0.94%  "
Set-CellText $ws.Range("D23") "2.05"
Set-CellText $ws.Range("E23") "  +7.19%  "
Set-CellText $ws.Range("D24") "72.24"
Set-CellText $ws.Range("E24") "  -0.06%  "
Set-CellText $ws.Range("E25") "  -0.08%  "
Set-CellText $ws.Range("D26") "4.34"
Set-CellText $ws.Range("E26") "  -3.07%  "
Set-CellText $ws.Range("D27") "9.83"
Set-CellText $ws.Range("E27") "  +1.80%  "
Set-CellText $ws.Range("D28") "2.827.36"
Set-CellText $ws.Range("E28") "  +2.06%  "
Set-CellText $ws.Range("E29") "  +0.04%  "
Set-CellText $ws.Range("D30") "0.0₃0942"
Set-CellText $ws.Range("E30") "  -2.54%  "
Set-CellText $ws.Range("D31") "8.08"
Set-CellText $ws.Range("E31") "  +0.24%  "
Set-CellText $ws.Range("D32") "512.81"
Set-CellText $ws.Range("E32") "  -5.50%  "
Set-CellText $ws.Range("E33") "  -2.04%  "
Set-CellText $ws.Range("E34") "  -1.06%  "
Set-CellText $ws.Range("D35") "1.00"
Set-CellText $ws.Range("E35") "  -0.01%  "
Set-CellText $ws.Range("D36") "162.65"
Set-CellText $ws.Range("E36") "  -1.91%  "
Set-CellText $ws.Range("D37") "19.60"
Set-CellText $ws.Range("E37") "  +1.82%  "
Set-CellText $ws.Range("E38") "  +0.46%  "
Set-CellText $ws.Range("E39") "  -0.68%  "
Set-CellText $ws.Range("E40") "  -3.79%  "
Set-CellText $ws.Range("D41") "0.107"
Set-CellText $ws.Range("E41") "  -9.35%  "
Set-CellText $ws.Range("E42") "  +0.07%  "
Set-CellText $ws.Range("D43") "5.02"
Set-CellText $ws.Range("E43") "  -0.66%  "
Set-CellText $ws.Range("E44") "  -3.07%  "
Set-CellText $ws.Range("E45") "  +0.52%  "
Set-CellText $ws.Range("D46") "39.29"
Set-CellText $ws.Range("E46") "  -1.75%  "
Set-CellText $ws.Range("D47") "153.86"
Set-CellText $ws.Range("E47") "  -0.35%  "
Set-CellText $ws.Range("D48") "3.73"
Set-CellText $ws.Range("E48") "  +2.46%  "
Set-CellText $ws.Range("D49") "0.553"
Set-CellText $ws.Range("E49") "  +3.54%  "
Set-CellText $ws.Range("E50") "  +1.75%  "
Set-CellText $ws.Range("E51") "  +1.46%  "
